# Weekly update: insert a new price record for Jengibre (Terminal La Palmera
# de La Serena) at row 132, pushing the existing history rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 132:161 down to 133:162 (mirrors the prior week's sheet having
# grown by one new weekly observation at the top of this block).
$ws.Rows("132:132").Insert()

# Populate the newly inserted row 132 with this week's observation.
$ws.Cells.Item(132, 1).Value = 8
$ws.Cells.Item(132, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(132, 3).Value = "Coquimbo"
$ws.Cells.Item(132, 4).Value = [DateTime]"2023-07-28"
$ws.Cells.Item(132, 4).NumberFormat = $ws.Cells.Item(133, 4).NumberFormat
$ws.Cells.Item(132, 5).Value = 4
$ws.Cells.Item(132, 6).Value = 100114007
$ws.Cells.Item(132, 7).Value = "Jengibre"
$ws.Cells.Item(132, 8).Value = "Sin especificar"
$ws.Cells.Item(132, 9).Value = "Primera"
$ws.Cells.Item(132, 10).Value = 340
$ws.Cells.Item(132, 11).Value = 17000
$ws.Cells.Item(132, 12).Value = 18000
$ws.Cells.Item(132, 13).Value = 17500
$ws.Cells.Item(132, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(132, 15).Value = "Perú"
$ws.Cells.Item(132, 16).Value = 1346
$ws.Cells.Item(132, 17).Value = 13
$ws.Cells.Item(132, 18).Value = "Hortaliza"
